$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab to reflect the new "through" date
$ws.Name = "Through 2022-09-22"

# Update the header label in I1 to match the new "through" date
$ws.Range("I1").Value = "2022 (through 09-22)"

# Update September (row 10) 2022 total and the overall Total row (row 14)
$ws.Range("I10").Value = 107
$ws.Range("I14").Value = 1242
